# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.739.40'
$ws.Range("E2").Value = '  +1.34%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.133.73'
$ws.Range("E3").Value = '  +1.18%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.32%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.83'
$ws.Range("E5").Value = '  +3.97%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '623.71'
$ws.Range("E6").Value = '  -0.33%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.16'
$ws.Range("E7").Value = '  +7.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.373'
$ws.Range("E8").Value = '  +2.37%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.22%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.118.47'
$ws.Range("E10").Value = '  +0.81%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.751'
$ws.Range("E11").Value = '  +3.22%  '

# Row 12
$ws.Range("E12").Value = '  +4.32%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  -1.91%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.61'
$ws.Range("E14").Value = '  -2.88%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.53'
$ws.Range("E15").Value = '  +1.07%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.969.88'
$ws.Range("E16").Value = '  +0.99%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.682.79'
$ws.Range("E17").Value = '  +0.83%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.092.10'
$ws.Range("E18").Value = '  +0.34%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.76'
$ws.Range("E19").Value = '  -0.73%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.00'
$ws.Range("E20").Value = '  +6.99%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000212'
$ws.Range("E21").Value = '  -2.45%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.85'
$ws.Range("E22").Value = '  +5.19%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '444.93'
$ws.Range("E23").Value = '  +1.85%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.24'
$ws.Range("E24").Value = '  +3.71%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.70'
$ws.Range("E25").Value = '  -4.28%  '

# Row 26
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.00'
$ws.Range("E26").Value = '  -1.87%  '

# Row 27
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '83.27'
$ws.Range("E27").Value = '  -6.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.246.24'
$ws.Range("E28").Value = '  -0.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.249'
$ws.Range("E30").Value = '  +27.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.176'
$ws.Range("E31").Value = '  +10.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.121'
$ws.Range("E32").Value = '  +39.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.34'
$ws.Range("E33").Value = '  -1.30%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.171'
$ws.Range("E34").Value = '  +11.02%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.06'
$ws.Range("E35").Value = '  +12.69%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.34'
$ws.Range("E36").Value = '  +26.76%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.70'
$ws.Range("E37").Value = '  +3.28%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '497.47'
$ws.Range("E38").Value = '  -2.06%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.91'
$ws.Range("E39").Value = '  +0.31%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  +1.89%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.56'
$ws.Range("E41").Value = '  -8.29%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.425'
$ws.Range("E42").Value = '  +3.42%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.22'
$ws.Range("E43").Value = '  +0.12%  '

# Row 44
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.94'
$ws.Range("E45").Value = '  +1.86%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.704'
$ws.Range("E46").Value = '  +2.18%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.05'
$ws.Range("E47").Value = '  +1.84%  '

# Row 48
$ws.Range("E48").Value = '  +1.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.636'
$ws.Range("E49").Value = '  -28.64%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.44'
$ws.Range("E50").Value = '  -1.12%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.45'
$ws.Range("E51").Value = '  +1.47%  '
